$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.829.81"
$ws.Range("E2").Value = "  +4.21%  "
$ws.Range("D3").Value = "2.268.87"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.80"
$ws.Range("E5").Value = "  +4.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.18"
$ws.Range("E6").Value = "  +4.92%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.82"
$ws.Range("E10").Value = "  +6.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.95"
$ws.Range("E11").Value = "  +6.18%  "
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").Value = "2.621.11"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").Value = "2.278.47"
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.765"
$ws.Range("E18").Value = "  +3.88%  "
$ws.Range("D19").Value = "41.762.17"
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.21"
$ws.Range("E20").Value = "  +8.32%  "
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.93"
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.93"
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.75"
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("E25").Value = "  +4.95%  "
$ws.Range("E27").Value = "  +5.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.21"
$ws.Range("E28").Value = "  +4.44%  "
$ws.Range("E29").Value = "  +11.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.63"
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.28"
$ws.Range("E31").Value = "  +7.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.90"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0748"
$ws.Range("E35").Value = "  +4.58%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.17"
$ws.Range("E37").Value = "  +9.50%  "
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("E40").Value = "  +4.36%  "
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("E42").Value = "  +4.61%  "
$ws.Range("D43").Value = "2.067.52"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.44"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("E45").Value = "  +3.11%  "
$ws.Range("E47").Value = "  +5.65%  "
$ws.Range("E48").Value = "  +6.83%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.53"
$ws.Range("E49").Value = "  +4.14%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.63"
$ws.Range("E51").Value = "  +7.01%  "
